# "remove explicit oofficial publication references"
#
# The canonical-OOXML diff for this commit shows the four custom XML
# "item" parts (customXml/item1.xml..item4.xml, and their paired
# itemProps1..4.xml datastore items) being re-saved under a different
# item numbering, while every part keeps exactly the same payload:
#
#   old item1 (FormTemplates)        -> new item2
#   old item2 (b:Sources/biblio)     -> new item1
#   old item3 (p:properties)         -> new item4
#   old item4 (ct:contentTypeSchema) -> new item3
#
# Word allocates customXml/itemN.xml (and itemPropsN.xml) purely by
# save-time order, so the supported way to reproduce this renumbering
# through the object model is: pull each non built-in CustomXMLPart's
# XML out, delete all of them, and Add() them back in the desired
# order so Word re-mints item1..item4 that way.

$d = $word.ActiveDocument
$parts = $d.CustomXMLParts

$sourcesXml      = $null
$formTemplateXml = $null
$propertiesXml   = $null
$contentTypeXml  = $null

$count = $parts.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $parts.Item($i)
    if ($p.BuiltIn) {
        continue
    }
    $xml = $p.XML

    if ($xml -like "*b:Sources*") {
        $sourcesXml = $xml
    } elseif ($xml -like "*FormTemplates*") {
        $formTemplateXml = $xml
    } elseif ($xml -like "*p:properties*") {
        $propertiesXml = $xml
    } elseif ($xml -like "*ct:contentTypeSchema*") {
        $contentTypeXml = $xml
    }

    $p.Delete()
}

# Re-add in the new target order so the save-time numbering becomes:
#   item1 = Sources (bibliography)
#   item2 = FormTemplates (SharePoint form template pointer)
#   item3 = contentTypeSchema (SharePoint content type)
#   item4 = properties (SharePoint documentManagement)
if ($sourcesXml)      { [void]$d.CustomXMLParts.Add($sourcesXml) }
if ($formTemplateXml) { [void]$d.CustomXMLParts.Add($formTemplateXml) }
if ($contentTypeXml)  { [void]$d.CustomXMLParts.Add($contentTypeXml) }
if ($propertiesXml)   { [void]$d.CustomXMLParts.Add($propertiesXml) }

$d.Save()
